# Auto-generated script applying scheduled market-price refresh to Ragnarok_Profits workbook
$wb = $excel.ActiveWorkbook

$ALC = $wb.Worksheets.Item("ALC")
$ARM = $wb.Worksheets.Item("ARM")
$BSM = $wb.Worksheets.Item("BSM")
$CRP = $wb.Worksheets.Item("CRP")
$CUL = $wb.Worksheets.Item("CUL")
$GSM = $wb.Worksheets.Item("GSM")
$LTW = $wb.Worksheets.Item("LTW")
$WVR = $wb.Worksheets.Item("WVR")


# --- ALC ---
$ALC.Cells.Item(64, 8).Value = 9000
$ALC.Cells.Item(64, 9).Value = 6000
$ALC.Cells.Item(64, 11).Value = 6000
$ALC.Cells.Item(64, 13).Value = -5752
$ALC.Cells.Item(67, 8).Value = 9000
$ALC.Cells.Item(67, 9).Value = 6000
$ALC.Cells.Item(67, 11).Value = 6000
$ALC.Cells.Item(67, 13).Value = -5142
$ALC.Cells.Item(96, 8).Value = 727075.4
$ALC.Cells.Item(96, 9).Value = 1704.3846
$ALC.Cells.Item(96, 11).Value = 5113.1538
$ALC.Cells.Item(96, 13).Value = -3740.1538
$ALC.Cells.Item(137, 8).Value = 1303721.5
$ALC.Cells.Item(137, 9).Value = 1050
$ALC.Cells.Item(137, 10).Value = 2172169.2
$ALC.Cells.Item(137, 11).Value = 3150
$ALC.Cells.Item(137, 12).Value = 6516507.600000001
$ALC.Cells.Item(137, 13).Value = -600
$ALC.Cells.Item(137, 14).Value = -6521607.600000001

# --- ARM ---
$ARM.Cells.Item(32, 8).Value = 3161
$ARM.Cells.Item(32, 9).Value = 3111.6875
$ARM.Cells.Item(32, 11).Value = 3111.6875
$ARM.Cells.Item(32, 13).Value = -2824.6875
$ARM.Cells.Item(63, 8).Value = 4999.5
$ARM.Cells.Item(63, 9).Value = 4999.5
$ARM.Cells.Item(63, 11).Value = 4999.5
$ARM.Cells.Item(63, 13).Value = -4313.5
$ARM.Cells.Item(66, 8).Value = 4999.5
$ARM.Cells.Item(66, 9).Value = 4999.5
$ARM.Cells.Item(66, 11).Value = 24997.5
$ARM.Cells.Item(66, 13).Value = -21565.5
$ARM.Cells.Item(74, 9).Value = 1042952.9
$ARM.Cells.Item(74, 11).Value = 1042952.9
$ARM.Cells.Item(74, 13).Value = -1042078.9
$ARM.Cells.Item(77, 9).Value = 1042952.9
$ARM.Cells.Item(77, 11).Value = 5214764.5
$ARM.Cells.Item(77, 13).Value = -5210396.5
$ARM.Cells.Item(122, 8).Value = 3913.2
$ARM.Cells.Item(122, 9).Value = 3792.4443
$ARM.Cells.Item(122, 11).Value = 11377.3329
$ARM.Cells.Item(122, 13).Value = -8927.332900000001
$ARM.Cells.Item(132, 8).Value = 5888451
$ARM.Cells.Item(132, 9).Value = 6690.2144
$ARM.Cells.Item(132, 10).Value = 33336666
$ARM.Cells.Item(132, 11).Value = 20070.6432
$ARM.Cells.Item(132, 12).Value = 100009998
$ARM.Cells.Item(132, 13).Value = -17540.6432
$ARM.Cells.Item(132, 14).Value = -100015058

# --- BSM ---
$BSM.Cells.Item(11, 8).Value = 3499.9167
$BSM.Cells.Item(11, 9).Value = 800
$BSM.Cells.Item(11, 10).Value = 4039.9
$BSM.Cells.Item(11, 11).Value = 800
$BSM.Cells.Item(11, 12).Value = 4039.9
$BSM.Cells.Item(11, 13).Value = -660
$BSM.Cells.Item(11, 14).Value = -4319.9
$BSM.Cells.Item(82, 8).Value = 4664.3335
$BSM.Cells.Item(82, 9).Value = 4664.3335
$BSM.Cells.Item(82, 10).Value = 0
$BSM.Cells.Item(82, 11).Value = 4664.3335
$BSM.Cells.Item(82, 12).Value = 0
$BSM.Cells.Item(82, 13).Value = -4281.3335
$BSM.Cells.Item(82, 14).ClearContents()
$BSM.Cells.Item(85, 8).Value = 4664.3335
$BSM.Cells.Item(85, 9).Value = 4664.3335
$BSM.Cells.Item(85, 10).Value = 0
$BSM.Cells.Item(85, 11).Value = 4664.3335
$BSM.Cells.Item(85, 12).Value = 0
$BSM.Cells.Item(85, 13).Value = -3338.3335
$BSM.Cells.Item(85, 14).ClearContents()
$BSM.Cells.Item(94, 8).Value = 1765.5778
$BSM.Cells.Item(94, 9).Value = 1439.5588
$BSM.Cells.Item(94, 11).Value = 1439.5588
$BSM.Cells.Item(94, 13).Value = -988.5588
$BSM.Cells.Item(105, 8).Value = 1274272.5
$BSM.Cells.Item(105, 9).Value = 1635571.9
$BSM.Cells.Item(105, 10).Value = 9724.75
$BSM.Cells.Item(105, 11).Value = 1635571.9
$BSM.Cells.Item(105, 12).Value = 9724.75
$BSM.Cells.Item(105, 13).Value = -1633824.9
$BSM.Cells.Item(105, 14).Value = -13218.75
$BSM.Cells.Item(134, 8).Value = 11113626
$BSM.Cells.Item(134, 9).Value = 2439.1667
$BSM.Cells.Item(134, 11).Value = 7317.500100000001
$BSM.Cells.Item(134, 13).Value = -4782.500100000001

# --- CRP ---
$CRP.Cells.Item(94, 9).Value = 1184.4286
$CRP.Cells.Item(94, 10).Value = 1172.4445
$CRP.Cells.Item(94, 11).Value = 1184.4286
$CRP.Cells.Item(94, 12).Value = 1172.4445
$CRP.Cells.Item(94, 13).Value = -733.4286
$CRP.Cells.Item(94, 14).Value = -2074.4445
$CRP.Cells.Item(105, 8).Value = 2120.6365
$CRP.Cells.Item(105, 9).Value = 1041.25
$CRP.Cells.Item(105, 11).Value = 1041.25
$CRP.Cells.Item(105, 13).Value = 705.75
$CRP.Cells.Item(107, 8).Value = 2114.9756
$CRP.Cells.Item(107, 9).Value = 1864.2333
$CRP.Cells.Item(107, 10).Value = 2798.818
$CRP.Cells.Item(107, 11).Value = 1864.2333
$CRP.Cells.Item(107, 12).Value = 2798.818
$CRP.Cells.Item(107, 13).Value = 55.7666999999999
$CRP.Cells.Item(107, 14).Value = -6638.818
$CRP.Cells.Item(132, 8).Value = 2913
$CRP.Cells.Item(132, 9).Value = 2681.1
$CRP.Cells.Item(132, 11).Value = 8043.299999999999
$CRP.Cells.Item(132, 13).Value = -5513.299999999999
$CRP.Cells.Item(134, 8).Value = 7006.143
$CRP.Cells.Item(134, 9).Value = 7006.143
$CRP.Cells.Item(134, 11).Value = 21018.429
$CRP.Cells.Item(134, 13).Value = -18483.429

# --- CUL ---
$CUL.Cells.Item(12, 8).Value = 2303.8462
$CUL.Cells.Item(12, 9).Value = 467.66666
$CUL.Cells.Item(12, 10).Value = 2854.7
$CUL.Cells.Item(12, 11).Value = 1402.99998
$CUL.Cells.Item(12, 12).Value = 8564.099999999999
$CUL.Cells.Item(12, 13).Value = -1229.99998
$CUL.Cells.Item(12, 14).Value = -8910.099999999999
$CUL.Cells.Item(87, 8).Value = 33333
$CUL.Cells.Item(87, 9).Value = 0
$CUL.Cells.Item(87, 11).Value = 0
$CUL.Cells.Item(87, 13).ClearContents()
$CUL.Cells.Item(90, 8).Value = 33333
$CUL.Cells.Item(90, 9).Value = 0
$CUL.Cells.Item(90, 11).Value = 0
$CUL.Cells.Item(90, 13).ClearContents()

# --- GSM ---
$GSM.Cells.Item(2, 8).Value = 63.333332
$GSM.Cells.Item(2, 9).Value = 65
$GSM.Cells.Item(2, 11).Value = 65
$GSM.Cells.Item(2, 13).Value = 48
$GSM.Cells.Item(97, 8).Value = 611.4286
$GSM.Cells.Item(97, 10).Value = 635.4286
$GSM.Cells.Item(97, 12).Value = 635.4286
$GSM.Cells.Item(97, 14).Value = -1627.4286
$GSM.Cells.Item(107, 8).Value = 2091.5
$GSM.Cells.Item(107, 9).Value = 1987.5
$GSM.Cells.Item(107, 10).Value = 2299.5
$GSM.Cells.Item(107, 11).Value = 1987.5
$GSM.Cells.Item(107, 12).Value = 2299.5
$GSM.Cells.Item(107, 13).Value = -67.5
$GSM.Cells.Item(107, 14).Value = -6139.5

# --- LTW ---
$LTW.Cells.Item(16, 8).Value = 4204.857
$LTW.Cells.Item(16, 9).Value = 1905.6666
$LTW.Cells.Item(16, 11).Value = 1905.6666
$LTW.Cells.Item(16, 13).Value = -1735.6666
$LTW.Cells.Item(22, 8).Value = 4242.222
$LTW.Cells.Item(22, 9).Value = 4047
$LTW.Cells.Item(22, 10).Value = 4398.4
$LTW.Cells.Item(22, 11).Value = 4047
$LTW.Cells.Item(22, 12).Value = 4398.4
$LTW.Cells.Item(22, 13).Value = -3752
$LTW.Cells.Item(22, 14).Value = -4988.4
$LTW.Cells.Item(27, 8).Value = 4242.222
$LTW.Cells.Item(27, 9).Value = 4047
$LTW.Cells.Item(27, 10).Value = 4398.4
$LTW.Cells.Item(27, 11).Value = 4047
$LTW.Cells.Item(27, 12).Value = 4398.4
$LTW.Cells.Item(27, 13).Value = -3940
$LTW.Cells.Item(27, 14).Value = -4612.4
$LTW.Cells.Item(40, 8).Value = 3255.6956
$LTW.Cells.Item(40, 9).Value = 3014.1
$LTW.Cells.Item(40, 10).Value = 4866.3335
$LTW.Cells.Item(40, 11).Value = 3014.1
$LTW.Cells.Item(40, 12).Value = 4866.3335
$LTW.Cells.Item(40, 13).Value = -2878.1
$LTW.Cells.Item(40, 14).Value = -5138.3335
$LTW.Cells.Item(61, 8).Value = 6030.6924
$LTW.Cells.Item(61, 10).Value = 8500.25
$LTW.Cells.Item(61, 12).Value = 8500.25
$LTW.Cells.Item(61, 14).Value = -8904.25
$LTW.Cells.Item(93, 8).Value = 9269509
$LTW.Cells.Item(93, 9).Value = 2500
$LTW.Cells.Item(93, 10).Value = 11122911
$LTW.Cells.Item(93, 11).Value = 2500
$LTW.Cells.Item(93, 12).Value = 11122911
$LTW.Cells.Item(93, 13).Value = -1252
$LTW.Cells.Item(93, 14).Value = -11125407
$LTW.Cells.Item(100, 8).Value = 19253940
$LTW.Cells.Item(100, 10).Value = 35753910
$LTW.Cells.Item(100, 12).Value = 35753910
$LTW.Cells.Item(100, 14).Value = -35754992
$LTW.Cells.Item(113, 8).Value = 6030.6924
$LTW.Cells.Item(113, 10).Value = 8500.25
$LTW.Cells.Item(113, 12).Value = 8500.25
$LTW.Cells.Item(113, 14).Value = -12840.25
$LTW.Cells.Item(122, 8).Value = 3442.3403
$LTW.Cells.Item(122, 9).Value = 3297.578
$LTW.Cells.Item(122, 10).Value = 6699.5
$LTW.Cells.Item(122, 11).Value = 9892.734
$LTW.Cells.Item(122, 12).Value = 20098.5
$LTW.Cells.Item(122, 13).Value = -7442.734
$LTW.Cells.Item(122, 14).Value = -24998.5
$LTW.Cells.Item(132, 8).Value = 6532.6
$LTW.Cells.Item(132, 10).Value = 7332.75
$LTW.Cells.Item(132, 12).Value = 21998.25
$LTW.Cells.Item(132, 14).Value = -27058.25

# --- WVR ---
$WVR.Cells.Item(17, 8).Value = 12581.5
$WVR.Cells.Item(17, 9).Value = 11497.8
$WVR.Cells.Item(17, 11).Value = 11497.8
$WVR.Cells.Item(17, 13).Value = -11325.8
$WVR.Cells.Item(59, 8).Value = 30000
$WVR.Cells.Item(59, 10).Value = 30000
$WVR.Cells.Item(59, 12).Value = 30000
$WVR.Cells.Item(59, 14).Value = -31476
$WVR.Cells.Item(107, 8).Value = 3361.244
$WVR.Cells.Item(107, 9).Value = 1718
$WVR.Cells.Item(107, 11).Value = 5154
$WVR.Cells.Item(107, 13).Value = -3234
$WVR.Cells.Item(132, 8).Value = 2502499.8
$WVR.Cells.Item(132, 9).Value = 2499.5
$WVR.Cells.Item(132, 10).Value = 5002500
$WVR.Cells.Item(132, 11).Value = 7498.5
$WVR.Cells.Item(132, 12).Value = 15007500
$WVR.Cells.Item(132, 13).Value = -4968.5
$WVR.Cells.Item(132, 14).Value = -15012560
$WVR.Cells.Item(136, 8).Value = 528476.6
$WVR.Cells.Item(136, 9).Value = 2238.2942
$WVR.Cells.Item(136, 11).Value = 6714.882599999999
$WVR.Cells.Item(136, 13).Value = -4164.882599999999
